$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "51.692.92"
Set-TextValue $ws.Range("E2") "  -0.85%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.924.29"
Set-TextValue $ws.Range("E3") "  +0.71%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.26%  "

# Row 5
Set-TextValue $ws.Range("D5") "351.82"
Set-TextValue $ws.Range("E5") "  -0.57%  "

# Row 6
Set-TextValue $ws.Range("D6") "106.20"
Set-TextValue $ws.Range("E6") "  -6.62%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.555"
Set-TextValue $ws.Range("E7") "  -0.21%  "

# Row 8
Set-TextValue $ws.Range("E8") "  -0.05%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.610"
Set-TextValue $ws.Range("E9") "  -2.10%  "

# Row 10
Set-TextValue $ws.Range("D10") "37.64"
Set-TextValue $ws.Range("E10") "  -4.94%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +1.00%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0852"
Set-TextValue $ws.Range("E12") "  -2.52%  "

# Row 13
Set-TextValue $ws.Range("D13") "18.98"
Set-TextValue $ws.Range("E13") "  -4.05%  "

# Row 14
Set-TextValue $ws.Range("D14") "3.385.06"
Set-TextValue $ws.Range("E14") "  +0.55%  "

# Row 15
Set-TextValue $ws.Range("D15") "7.55"
Set-TextValue $ws.Range("E15") "  -2.29%  "

# Row 16
Set-TextValue $ws.Range("D16") "2.917.06"
Set-TextValue $ws.Range("E16") "  +0.27%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.964"
Set-TextValue $ws.Range("E17") "  -2.01%  "

# Row 18
Set-TextValue $ws.Range("D18") "51.604.32"
Set-TextValue $ws.Range("E18") "  -1.18%  "

# Row 19
Set-TextValue $ws.Range("D19") "3.41"
Set-TextValue $ws.Range("E19") "  +2.29%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -3.25%  "

# Row 21
Set-TextValue $ws.Range("D21") "13.38"
Set-TextValue $ws.Range("E21") "  -4.88%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.0₃0958"
Set-TextValue $ws.Range("E22") "  -2.16%  "

# Row 23
Set-TextValue $ws.Range("D23") "68.96"
Set-TextValue $ws.Range("E23") "  -2.89%  "

# Row 24
Set-TextValue $ws.Range("D24") "261.49"
Set-TextValue $ws.Range("E24") "  -2.89%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.73"
Set-TextValue $ws.Range("E25") "  -2.97%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.172"
Set-TextValue $ws.Range("E26") "  -5.36%  "

# Row 27
Set-TextValue $ws.Range("D27") "26.46"
Set-TextValue $ws.Range("E27") "  -1.04%  "

# Row 28
Set-TextValue $ws.Range("E28") "  +0.12%  "

# Row 29
Set-TextValue $ws.Range("D29") "7.39"
Set-TextValue $ws.Range("E29") "  +8.76%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -0.41%  "

# Row 31
Set-TextValue $ws.Range("D31") "10.21"
Set-TextValue $ws.Range("E31") "  -4.08%  "

# Row 32
Set-TextValue $ws.Range("B32") "InjectiveProtocol"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D32") "35.72"
Set-TextValue $ws.Range("E32") "  -4.54%  "

# Row 33
Set-TextValue $ws.Range("B33") "Toncoin"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D33") "2.16"
Set-TextValue $ws.Range("E33") "  -4.95%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.94"
Set-TextValue $ws.Range("E34") "  -2.56%  "

# Row 35
Set-TextValue $ws.Range("D35") "50.89"
Set-TextValue $ws.Range("E35") "  -4.01%  "

# Row 36
Set-TextValue $ws.Range("B36") "VeChain"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D36") "0.0426"
Set-TextValue $ws.Range("E36") "  -5.52%  "

# Row 37
Set-TextValue $ws.Range("B37") "FirstDigitalUSD"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D37") "0.997"
Set-TextValue $ws.Range("E37") "  -0.18%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.14"
Set-TextValue $ws.Range("E38") "  -5.17%  "

# Row 39
Set-TextValue $ws.Range("B39") "Celestia"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D39") "17.72"
Set-TextValue $ws.Range("E39") "  -5.74%  "

# Row 40
Set-TextValue $ws.Range("B40") "ARBITRUM"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D40") "1.96"
Set-TextValue $ws.Range("E40") "  -3.85%  "

# Row 41
Set-TextValue $ws.Range("D41") "2.65"
Set-TextValue $ws.Range("E41") "  -2.83%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -1.04%  "

# Row 43
Set-TextValue $ws.Range("D43") "22.91"
Set-TextValue $ws.Range("E43") "  -0.76%  "

# Row 44
Set-TextValue $ws.Range("D44") "119.69"
Set-TextValue $ws.Range("E44") "  +1.56%  "

# Row 45
Set-TextValue $ws.Range("E45") "  -1.14%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.44"
Set-TextValue $ws.Range("E46") "  -3.51%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.098.27"
Set-TextValue $ws.Range("E47") "  -3.78%  "

# Row 48
Set-TextValue $ws.Range("D48") "3.31"
Set-TextValue $ws.Range("E48") "  -5.78%  "

# Row 49
Set-TextValue $ws.Range("D49") "3.210.47"
Set-TextValue $ws.Range("E49") "  +0.43%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.237"
Set-TextValue $ws.Range("E50") "  -8.48%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0338"
Set-TextValue $ws.Range("E51") "  -4.71%  "
